# Add a "Reference" sheet (citation for the dataset) ahead of the data
# sheet, rename the data sheet to "Data", and drop the two unused blank
# sheets ("Sheet2" / "Sheet3").

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the unused empty sheets -----------------------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
[void]$sheet2.Delete()
$sheet3 = $wb.Worksheets.Item("Sheet3")
[void]$sheet3.Delete()

# --- Rename the data sheet ---------------------------------------------
$data = $wb.Worksheets.Item("Sheet1")
$data.Name = "Data"

# Reset the data sheet's view (clear the old scroll position / selection)
$data.Activate()
[void]$data.Range("AX6").Select()

# --- Add the new "Reference" sheet, placed before "Data" ---------------
$ref = $wb.Worksheets.Add()
$ref.Name = "Reference"

$ref.Columns("A").ColumnWidth = 102

# Row 1 (merged A1:A2): intro line
$ref.Range("A1:A2").Merge()
$ref.Range("A1").Value = "If you use this dataset in a publication, please cite:  "
$ref.Range("A1:A2").Font.Size = 16
$ref.Range("A1:A2").HorizontalAlignment = -4131
$ref.Range("A1:A2").VerticalAlignment = -4108
$ref.Range("A1:A2").WrapText = $true
$ref.Range("A1:A2").Borders.LineStyle = -4142

# Row 3: citation, bold title followed by regular authors/journal text
$titleText = "The use of psoriasis biomarkers, including trajectory of clinical response, to predict clearance and remission duration to UVB phototherapy"
$restText = ". N. Watson, N. Wilson, F. Shmarov, P. Zuliani, N. J. Reynolds, S. C. Weatherhead. Journal of the European Academy of Dermatology & Venereology 35: 2250-2258, 2021."
$ref.Range("A3").Value = $titleText + $restText
$ref.Range("A3").Font.Size = 16
$ref.Range("A3").Characters(1, $titleText.Length).Font.Bold = $true
$ref.Range("A3").WrapText = $true
$ref.Rows(3).RowHeight = 88

# Row 4: hyperlink to the DOI
$url = "https://doi.org/10.1111/jdv.17519"
$ref.Range("A4").Value = $url
[void]$ref.Hyperlinks.Add($ref.Range("A4"), $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url)
$ref.Range("A4").Font.Size = 16
$ref.Rows(4).RowHeight = 21

$ref.Activate()
$excel.ActiveWindow.Zoom = 130
[void]$ref.Range("A1").Select()
